# Refresh cryptos list with latest scraped values (GitHub Actions run).
# D-column prices are stored as literal text (e.g. "97.098.03", "1.00") in the
# source sheet, so numeric-looking values are written with a leading quote
# (standard Excel "quote prefix") to force text entry instead of silently
# being parsed into a Number/losing formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range('D2').Value = '''97.098.03'
$ws.Range('E2').Value = '  +0.44%  '

# Row 3 - Ethereum
$ws.Range('D3').Value = '''3.739.26'
$ws.Range('E3').Value = '  +1.70%  '

# Row 4 - TetherUSD
$ws.Range('E4').Value = '  -0.03%  '

# Row 5 - Solana
$ws.Range('D5').Value = '''239.04'
$ws.Range('E5').Value = '  -0.09%  '

# Row 6 - XRP
$ws.Range('E6').Value = '  +1.88%  '

# Row 7 - BNB
$ws.Range('D7').Value = '''662.57'

# Row 8 - Dogecoin
$ws.Range('D8').Value = '''0.427'
$ws.Range('E8').Value = '  +1.47%  '

# Row 9 - Cardano
$ws.Range('E9').Value = '  -1.04%  '

# Row 10 - USDC
$ws.Range('E10').Value = '  -0.02%  '

# Row 11 - LidoStakedEther
$ws.Range('D11').Value = '''3.737.15'
$ws.Range('E11').Value = '  +1.65%  '

# Row 12 - ShibaInu
$ws.Range('D12').Value = '''0.0000323'
$ws.Range('E12').Value = '  +20.31%  '

# Row 13 - Avalanche
$ws.Range('D13').Value = '''44.86'
$ws.Range('E13').Value = '  -1.29%  '

# Row 14 - TRON
$ws.Range('E14').Value = '  +1.63%  '

# Row 15 - Toncoin
$ws.Range('D15').Value = '''6.93'
$ws.Range('E15').Value = '  +1.52%  '

# Row 16 - WrappedliquidstakedEther2.0
$ws.Range('D16').Value = '''4.433.92'
$ws.Range('E16').Value = '  +1.64%  '

# Row 17 - WrappedBTC
$ws.Range('D17').Value = '''96.936.36'
$ws.Range('E17').Value = '  +0.69%  '

# Row 18 - Polkadot
$ws.Range('E18').Value = '  +17.71%  '

# Row 19 - WrappedEther
$ws.Range('D19').Value = '''3.742.88'
$ws.Range('E19').Value = '  +2.59%  '

# Row 20 - Uniswap
$ws.Range('D20').Value = '''13.24'
$ws.Range('E20').Value = '  +3.71%  '

# Row 21 - Chainlink
$ws.Range('D21').Value = '''18.93'
$ws.Range('E21').Value = '  +0.31%  '

# Row 22 - Stellar
$ws.Range('E22').Value = '  -3.79%  '

# Row 23 - BitcoinCash
$ws.Range('D23').Value = '''529.62'
$ws.Range('E23').Value = '  -0.13%  '

# Row 24 - SuiNetwork
$ws.Range('D24').Value = '''3.49'
$ws.Range('E24').Value = '  +1.02%  '

# Row 25 - PEPE
$ws.Range('D25').Value = '''0.0000227'
$ws.Range('E25').Value = '  +11.30%  '

# Row 26 - NEARProtocol
$ws.Range('E26').Value = '  -2.76%  '

# Row 27 - Litecoin
$ws.Range('D27').Value = '''108.72'
$ws.Range('E27').Value = '  +6.87%  '

# Row 28 - Aptos
$ws.Range('B28').Value = 'Aptos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D28').Value = '''13.74'
$ws.Range('E28').Value = '  +4.00%  '

# Row 29 - Hedera
$ws.Range('B29').Value = 'Hedera'
$ws.Range('C29').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D29').Value = '''0.192'
$ws.Range('E29').Value = '  +14.40%  '

# Row 30 - WrappedeETH
$ws.Range('D30').Value = '''3.914.70'
$ws.Range('E30').Value = '  +1.12%  '

# Row 31 - InternetComputer(DFINITY)
$ws.Range('D31').Value = '''13.07'
$ws.Range('E31').Value = '  +4.81%  '

# Row 32 - PancakeSwap
$ws.Range('E32').Value = '  +0.93%  '

# Row 33 - Dai
$ws.Range('D33').Value = '''0.999'

# Row 34 - Cronos
$ws.Range('E34').Value = '  +4.25%  '

# Row 35 - Fetch.AI
$ws.Range('E35').Value = '  -2.67%  '

# Row 36 - EthereumClassic
$ws.Range('D36').Value = '''33.44'
$ws.Range('E36').Value = '  +2.85%  '

# Row 37 - Binance-PegBSC-USD
$ws.Range('D37').Value = '''0.997'
$ws.Range('E37').Value = '  -0.32%  '

# Row 38 - Bittensor
$ws.Range('D38').Value = '''652.92'
$ws.Range('E38').Value = '  -3.76%  '

# Row 39 - PolygonEcosystemToken
$ws.Range('D39').Value = '''0.599'
$ws.Range('E39').Value = '  +1.09%  '

# Row 40 - RenderToken
$ws.Range('D40').Value = '''8.89'
$ws.Range('E40').Value = '  +0.73%  '

# Row 41 - USDe
$ws.Range('E41').Value = '  +0.00%  '

# Row 42 - EnergySwap
$ws.Range('B42').Value = 'EnergySwap'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D42').Value = '''42.64'
$ws.Range('E42').Value = '  +7.31%  '

# Row 43 - Kaspa
$ws.Range('B43').Value = 'Kaspa'
$ws.Range('C43').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D43').Value = '''0.168'
$ws.Range('E43').Value = '  +4.98%  '

# Row 44 - Filecoin
$ws.Range('B44').Value = 'Filecoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D44').Value = '''6.88'
$ws.Range('E44').Value = '  +4.33%  '

# Row 45 - ImmutableX
$ws.Range('B45').Value = 'ImmutableX'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D45').Value = '''2.05'
$ws.Range('E45').Value = '  +2.73%  '

# Row 46 - ARBITRUM
$ws.Range('B46').Value = 'ARBITRUM'
$ws.Range('C46').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D46').Value = '''0.993'
$ws.Range('E46').Value = '  +3.39%  '

# Row 47 - Algorand
$ws.Range('D47').Value = '''0.476'
$ws.Range('E47').Value = '  +7.60%  '

# Row 48 - VeChain
$ws.Range('D48').Value = '''0.0461'
$ws.Range('E48').Value = '  -0.99%  '

# Row 49 - Stacks
$ws.Range('D49').Value = '''2.41'
$ws.Range('E49').Value = '  +3.78%  '

# Row 50 - Cosmos
$ws.Range('D50').Value = '''8.80'
$ws.Range('E50').Value = '  +2.16%  '

# Row 51 - WhiteBITCoin
$ws.Range('D51').Value = '''23.60'
$ws.Range('E51').Value = '  -0.26%  '
